$wb = $excel.ActiveWorkbook

# --- ALC (sheet1.xml) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3872.0557
$ws.Range("I51").Value = 2084.1667
$ws.Range("J51").Value = 4766
$ws.Range("K51").Value = 2084.1667
$ws.Range("L51").Value = 4766
$ws.Range("M51").Value = -1600.1667
$ws.Range("N51").Value = -5734
$ws.Range("H112").Value = 17859928
$ws.Range("I112").Value = 3891.6667
$ws.Range("J112").Value = 20002652
$ws.Range("K112").Value = 11675.0001
$ws.Range("L112").Value = 60007956
$ws.Range("M112").Value = -10567.0001
$ws.Range("N112").Value = -60010172
$ws.Range("H127").Value = 1123.1
$ws.Range("J127").Value = 2857.5
$ws.Range("L127").Value = 8572.5
$ws.Range("N127").Value = -18492.5
$ws.Range("H141").Value = 1498.4
$ws.Range("I141").Value = 1279.3636
$ws.Range("K141").Value = 3838.0908
$ws.Range("M141").Value = 1341.9092

# --- ARM (sheet2.xml) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2798.28
$ws.Range("I61").Value = 1563.4706
$ws.Range("J61").Value = 5422.25
$ws.Range("K61").Value = 1563.4706
$ws.Range("L61").Value = 5422.25
$ws.Range("M61").Value = -1351.4706
$ws.Range("N61").Value = -5846.25
$ws.Range("H102").Value = 9086.666999999999
$ws.Range("I102").Value = 8504
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 8504
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -6882
$ws.Range("N102").Value = -15244
$ws.Range("H128").Value = 115000
$ws.Range("J128").Value = 115000
$ws.Range("L128").Value = 115000
$ws.Range("N128").Value = -124960
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H130").Value = 82105.75
$ws.Range("J130").Value = 82105.75
$ws.Range("L130").Value = 82105.75
$ws.Range("N130").Value = -92145.75
$ws.Range("H132").Value = 34325.934
$ws.Range("I132").Value = 2504.0781
$ws.Range("K132").Value = 7512.2343
$ws.Range("M132").Value = -4982.2343
$ws.Range("H133").Value = 106993
$ws.Range("J133").Value = 106993
$ws.Range("L133").Value = 106993
$ws.Range("N133").Value = -112053
$ws.Range("H136").Value = 2798.28
$ws.Range("I136").Value = 1563.4706
$ws.Range("J136").Value = 5422.25
$ws.Range("K136").Value = 4690.4118
$ws.Range("L136").Value = 16266.75
$ws.Range("M136").Value = -2140.4118
$ws.Range("N136").Value = -21366.75

# --- BSM (sheet3.xml) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 48586
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H59").Value = 85000
$ws.Range("J59").Value = 85000
$ws.Range("L59").Value = 85000
$ws.Range("N59").Value = -86694
$ws.Range("H81").Value = 33926.668
$ws.Range("J81").Value = 33926.668
$ws.Range("L81").Value = 33926.668
$ws.Range("N81").Value = -36048.668
$ws.Range("H84").Value = 33926.668
$ws.Range("J84").Value = 33926.668
$ws.Range("L84").Value = 101780.004
$ws.Range("N84").Value = -112388.004
$ws.Range("H94").Value = 2303.2354
$ws.Range("I94").Value = 2085.5557
$ws.Range("K94").Value = 2085.5557
$ws.Range("M94").Value = -1634.5557
$ws.Range("H112").Value = 66249.75
$ws.Range("J112").Value = 56666.332
$ws.Range("L112").Value = 56666.332
$ws.Range("N112").Value = -59620.332
$ws.Range("H132").Value = 144537.61
$ws.Range("J132").Value = 144537.61
$ws.Range("L132").Value = 144537.61
$ws.Range("N132").Value = -154657.61
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 888.13635
$ws.Range("I134").Value = 776.25
$ws.Range("J134").Value = 2007
$ws.Range("K134").Value = 2328.75
$ws.Range("L134").Value = 6021
$ws.Range("M134").Value = 206.25
$ws.Range("N134").Value = -11091

# --- CRP (sheet4.xml) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4270.375
$ws.Range("I31").Value = 2567.3713
$ws.Range("K31").Value = 2567.3713
$ws.Range("M31").Value = -2272.3713
$ws.Range("H34").Value = 4270.375
$ws.Range("I34").Value = 2567.3713
$ws.Range("K34").Value = 2567.3713
$ws.Range("M34").Value = -2365.3713
$ws.Range("H52").Value = 89377.8
$ws.Range("J52").Value = 89377.8
$ws.Range("L52").Value = 89377.8
$ws.Range("N52").Value = -89965.8
$ws.Range("H58").Value = 11002587
$ws.Range("I58").Value = 2160.5
$ws.Range("J58").Value = 21156826
$ws.Range("K58").Value = 2160.5
$ws.Range("L58").Value = 21156826
$ws.Range("M58").Value = -1957.5
$ws.Range("N58").Value = -21157232
$ws.Range("H68").Value = 48795
$ws.Range("J68").Value = 48795
$ws.Range("L68").Value = 48795
$ws.Range("N68").Value = -50293
$ws.Range("H71").Value = 48795
$ws.Range("J71").Value = 48795
$ws.Range("L71").Value = 146385
$ws.Range("N71").Value = -153873
$ws.Range("H74").Value = 58986.75
$ws.Range("J74").Value = 58986.75
$ws.Range("L74").Value = 58986.75
$ws.Range("N74").Value = -60734.75
$ws.Range("H77").Value = 58986.75
$ws.Range("J77").Value = 58986.75
$ws.Range("L77").Value = 176960.25
$ws.Range("N77").Value = -185696.25
$ws.Range("H86").Value = 9877.111000000001
$ws.Range("J86").Value = 10919.2
$ws.Range("L86").Value = 10919.2
$ws.Range("N86").Value = -13165.2
$ws.Range("H89").Value = 9877.111000000001
$ws.Range("J89").Value = 10919.2
$ws.Range("L89").Value = 54596
$ws.Range("N89").Value = -65828
$ws.Range("H132").Value = 4014.6287
$ws.Range("I132").Value = 3550.25
$ws.Range("J132").Value = 5872.143
$ws.Range("K132").Value = 10650.75
$ws.Range("L132").Value = 17616.429
$ws.Range("M132").Value = -8120.75
$ws.Range("N132").Value = -22676.429
$ws.Range("H134").Value = 3324.15
$ws.Range("I134").Value = 2791.1538
$ws.Range("K134").Value = 8373.4614
$ws.Range("M134").Value = -5838.4614
$ws.Range("H135").Value = 93570
$ws.Range("J135").Value = 93570
$ws.Range("L135").Value = 93570
$ws.Range("N135").Value = -103710
$ws.Range("H136").Value = 11002587
$ws.Range("I136").Value = 2160.5
$ws.Range("J136").Value = 21156826
$ws.Range("K136").Value = 6481.5
$ws.Range("L136").Value = 63470478
$ws.Range("M136").Value = -3931.5
$ws.Range("N136").Value = -63475578
$ws.Range("H138").Value = 199994
$ws.Range("J138").Value = 199994
$ws.Range("L138").Value = 199994
$ws.Range("N138").Value = -210274
$ws.Range("H139").Value = 84994.664
$ws.Range("J139").Value = 84994.664
$ws.Range("L139").Value = 84994.664
$ws.Range("N139").Value = -95274.664

# --- CUL (sheet5.xml) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 7875.077
$ws.Range("I14").Value = 7875.077
$ws.Range("K14").Value = 23625.231
$ws.Range("M14").Value = -23452.231
$ws.Range("H70").Value = 252.75
$ws.Range("I70").Value = 252.75
$ws.Range("K70").Value = 758.25
$ws.Range("M70").Value = -443.25
$ws.Range("H73").Value = 252.75
$ws.Range("I73").Value = 252.75
$ws.Range("K73").Value = 758.25
$ws.Range("M73").Value = 333.75
$ws.Range("H104").Value = 192.25
$ws.Range("I104").Value = 192.25
$ws.Range("K104").Value = 576.75
$ws.Range("M104").Value = 2044.25
$ws.Range("H121").Value = 542.7143
$ws.Range("J121").Value = 639.8
$ws.Range("L121").Value = 1919.4
$ws.Range("N121").Value = -4539.4

# --- GSM (sheet6.xml) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 866.3333
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 866.3333
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 866.3333
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1098.3333
$ws.Range("H132").Value = 1560.4445
$ws.Range("I132").Value = 1340.6666
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4021.9998
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1491.9998
$ws.Range("N132").Value = -11060

# --- LTW (sheet7.xml) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3145.5833
$ws.Range("I16").Value = 2774.3
$ws.Range("K16").Value = 2774.3
$ws.Range("M16").Value = -2604.3
$ws.Range("H40").Value = 12825001
$ws.Range("I40").Value = 16669949
$ws.Range("J40").Value = 8508
$ws.Range("K40").Value = 16669949
$ws.Range("L40").Value = 8508
$ws.Range("M40").Value = -16669813
$ws.Range("N40").Value = -8780
$ws.Range("H93").Value = 2442.158
$ws.Range("I93").Value = 2586.8823
$ws.Range("J93").Value = 1212
$ws.Range("K93").Value = 2586.8823
$ws.Range("L93").Value = 1212
$ws.Range("M93").Value = -1338.8823
$ws.Range("N93").Value = -3708
$ws.Range("H122").Value = 3252.818
$ws.Range("I122").Value = 2831.3333
$ws.Range("K122").Value = 8493.999899999999
$ws.Range("M122").Value = -6043.999899999999
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840
$ws.Range("H136").Value = 2215.6553
$ws.Range("I136").Value = 2102
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 6306
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -3756
$ws.Range("N136").Value = -16350

# --- WVR (sheet8.xml) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 327.2
$ws.Range("I132").Value = 327.2
$ws.Range("K132").Value = 981.5999999999999
$ws.Range("M132").Value = 1548.4
$ws.Range("H136").Value = 4390.9814
$ws.Range("J136").Value = 13933.125
$ws.Range("L136").Value = 41799.375
$ws.Range("N136").Value = -46899.375
